$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("A2").Value = "MuSCs"
    $ws.Range("B2").Value = "Ifnb1"
    $ws.Range("C2").Value = "Ifnar2"
    $ws.Range("D2").Value = "ECs"
    $ws.Range("E2").Value = 1
    $ws.Range("F2").Value = 0.3333333333333333
    $ws.Range("G2").Value = 0.01317833333333333
    $ws.Range("H2").Value = 0.039535
    $ws.Range("I2").Value = 0.2048477437473122
    $ws.Range("J2").Value = 0.2048477437473122
    $ws.Range("K2").Value = 3
    $ws.Range("L2").Value = 1
    $ws.Range("M2").Value = 29.04796866666666
    $ws.Range("N2").Value = 87.14390599999999
    $ws.Range("O2").Value = 0.2371972210028098
    $ws.Range("P2").Value = 0.2371972210028099
    $ws.Range("Q2").Value = 0.3828038137455555
    $ws.Range("R2").Value = 3.445234323709999
    $ws.Range("S2").Value = 0.04858931554555816
    $ws.Range("T2").Value = 0.04858931554555817
    # Row 3
    $ws.Range("A3").Value = "MuSCs"
    $ws.Range("B3").Value = "Ifnb1"
    $ws.Range("C3").Value = "Ifnar2"
    $ws.Range("D3").Value = "FAPs"
    $ws.Range("E3").Value = 1
    $ws.Range("F3").Value = 0.3333333333333333
    $ws.Range("G3").Value = 0.01317833333333333
    $ws.Range("H3").Value = 0.039535
    $ws.Range("I3").Value = 0.2048477437473122
    $ws.Range("J3").Value = 0.2048477437473122
    $ws.Range("K3").Value = 3
    $ws.Range("L3").Value = 1
    $ws.Range("M3").Value = 31.81541733333333
    $ws.Range("N3").Value = 95.446252
    $ws.Range("O3").Value = 0.2597953978506987
    $ws.Range("P3").Value = 0.2597953978506987
    $ws.Range("Q3").Value = 0.4192741747577778
    $ws.Range("R3").Value = 3.77346757282
    $ws.Range("S3").Value = 0.05321850108565093
    $ws.Range("T3").Value = 0.05321850108565093
    # Row 4
    $ws.Range("A4").Value = "MuSCs"
    $ws.Range("B4").Value = "Ifnb1"
    $ws.Range("C4").Value = "Ifnar2"
    $ws.Range("D4").Value = "MuSCs"
    $ws.Range("E4").Value = 1
    $ws.Range("F4").Value = 0.3333333333333333
    $ws.Range("G4").Value = 0.01317833333333333
    $ws.Range("H4").Value = 0.039535
    $ws.Range("I4").Value = 0.2048477437473122
    $ws.Range("J4").Value = 0.2048477437473122
    $ws.Range("K4").Value = 3
    $ws.Range("L4").Value = 1
    $ws.Range("M4").Value = 12.70280433333333
    $ws.Range("N4").Value = 38.108413
    $ws.Range("O4").Value = 0.1037273869778955
    $ws.Range("P4").Value = 0.1037273869778955
    $ws.Range("Q4").Value = 0.1674017897727778
    $ws.Range("R4").Value = 1.506616107955
    $ws.Range("S4").Value = 0.02124832118722623
    $ws.Range("T4").Value = 0.02124832118722623
    # Row 5
    $ws.Range("A5").Value = "MuSCs"
    $ws.Range("B5").Value = "Ifnb1"
    $ws.Range("C5").Value = "Ifnar2"
    $ws.Range("D5").Value = "Resolving-Mac"
    $ws.Range("E5").Value = 1
    $ws.Range("F5").Value = 0.3333333333333333
    $ws.Range("G5").Value = 0.01317833333333333
    $ws.Range("H5").Value = 0.039535
    $ws.Range("I5").Value = 0.2048477437473122
    $ws.Range("J5").Value = 0.2048477437473122
    $ws.Range("K5").Value = 3
    $ws.Range("L5").Value = 1
    $ws.Range("M5").Value = 48.89716966666666
    $ws.Range("N5").Value = 146.691509
    $ws.Range("O5").Value = 0.3992799941685959
    $ws.Range("P5").Value = 0.399279994168596
    $ws.Range("Q5").Value = 0.6443832009238889
    $ws.Range("R5").Value = 5.799448808315
    $ws.Range("S5").Value = 0.08179160592887684
    $ws.Range("T5").Value = 0.08179160592887684
    # Row 6
    $ws.Range("A6").Value = "Resolving-Mac"
    $ws.Range("B6").Value = "Ifnb1"
    $ws.Range("C6").Value = "Ifnar2"
    $ws.Range("D6").Value = "ECs"
    $ws.Range("E6").Value = 2
    $ws.Range("F6").Value = 0.6666666666666666
    $ws.Range("G6").Value = 0.051154
    $ws.Range("H6").Value = 0.153462
    $ws.Range("I6").Value = 0.7951522562526879
    $ws.Range("J6").Value = 0.7951522562526879
    $ws.Range("K6").Value = 3
    $ws.Range("L6").Value = 1
    $ws.Range("M6").Value = 29.04796866666666
    $ws.Range("N6").Value = 87.14390599999999
    $ws.Range("O6").Value = 0.2371972210028098
    $ws.Range("P6").Value = 0.2371972210028099
    $ws.Range("Q6").Value = 1.485919789174666
    $ws.Range("R6").Value = 13.373278102572
    $ws.Range("S6").Value = 0.1886079054572517
    $ws.Range("T6").Value = 0.1886079054572517
    # Row 7
    $ws.Range("A7").Value = "Resolving-Mac"
    $ws.Range("B7").Value = "Ifnb1"
    $ws.Range("C7").Value = "Ifnar2"
    $ws.Range("D7").Value = "FAPs"
    $ws.Range("E7").Value = 2
    $ws.Range("F7").Value = 0.6666666666666666
    $ws.Range("G7").Value = 0.051154
    $ws.Range("H7").Value = 0.153462
    $ws.Range("I7").Value = 0.7951522562526879
    $ws.Range("J7").Value = 0.7951522562526879
    $ws.Range("K7").Value = 3
    $ws.Range("L7").Value = 1
    $ws.Range("M7").Value = 31.81541733333333
    $ws.Range("N7").Value = 95.446252
    $ws.Range("O7").Value = 0.2597953978506987
    $ws.Range("P7").Value = 0.2597953978506987
    $ws.Range("Q7").Value = 1.627485858269333
    $ws.Range("R7").Value = 14.647372724424
    $ws.Range("S7").Value = 0.2065768967650478
    $ws.Range("T7").Value = 0.2065768967650478
    # Row 8
    $ws.Range("A8").Value = "Resolving-Mac"
    $ws.Range("B8").Value = "Ifnb1"
    $ws.Range("C8").Value = "Ifnar2"
    $ws.Range("D8").Value = "MuSCs"
    $ws.Range("E8").Value = 2
    $ws.Range("F8").Value = 0.6666666666666666
    $ws.Range("G8").Value = 0.051154
    $ws.Range("H8").Value = 0.153462
    $ws.Range("I8").Value = 0.7951522562526879
    $ws.Range("J8").Value = 0.7951522562526879
    $ws.Range("K8").Value = 3
    $ws.Range("L8").Value = 1
    $ws.Range("M8").Value = 12.70280433333333
    $ws.Range("N8").Value = 38.108413
    $ws.Range("O8").Value = 0.1037273869778955
    $ws.Range("P8").Value = 0.1037273869778955
    $ws.Range("Q8").Value = 0.6497992528673333
    $ws.Range("R8").Value = 5.848193275806
    $ws.Range("S8").Value = 0.08247906579066931
    $ws.Range("T8").Value = 0.08247906579066933
    # Row 9
    $ws.Range("A9").Value = "Resolving-Mac"
    $ws.Range("B9").Value = "Ifnb1"
    $ws.Range("C9").Value = "Ifnar2"
    $ws.Range("D9").Value = "Resolving-Mac"
    $ws.Range("E9").Value = 2
    $ws.Range("F9").Value = 0.6666666666666666
    $ws.Range("G9").Value = 0.051154
    $ws.Range("H9").Value = 0.153462
    $ws.Range("I9").Value = 0.7951522562526879
    $ws.Range("J9").Value = 0.7951522562526879
    $ws.Range("K9").Value = 3
    $ws.Range("L9").Value = 1
    $ws.Range("M9").Value = 48.89716966666666
    $ws.Range("N9").Value = 146.691509
    $ws.Range("O9").Value = 0.3992799941685959
    $ws.Range("P9").Value = 0.399279994168596
    $ws.Range("Q9").Value = 2.501285817128666
    $ws.Range("R9").Value = 22.511572354158
    $ws.Range("S9").Value = 0.3174883882397191
    $ws.Range("T9").Value = 0.3174883882397192
